$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.855.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.705.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  +1.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.009"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3948"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4082"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.011"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08831"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.642"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001372"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.081"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.701.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.92"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07124"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.435"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.012"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("E23").Value = "  +0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.850.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.064"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.354"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.754"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.230"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.743"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.890.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08957"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.059"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.991"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09191"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7940"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.96%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.474"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7294"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.631"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.263"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.009"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.337"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.96%  "

$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2761"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.91"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.00%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02896"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.15%  "
